# Refresh the handback-status report's timestamp strings (regenerated report).
$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for cb642401-...md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-29 12:48:21"

# "zh-cn" sheet: Correspond Handoff/Handback datetimes for cb642401-...xlf
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-29 12:48:17"
$wsZhCn.Range("K4").Value = "2016-08-29 12:48:35"

# "de-de" sheet: Correspond Handoff/Handback datetimes for cb642401-...xlf
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-29 12:48:21"
$wsDeDe.Range("K4").Value = "2016-08-29 12:48:43"
